# Remove V47 validation scenario for now
#
# The "Vestas_V47_public" project (row 2 of the project list) is being
# dropped from the template. Deleting the entire worksheet row takes care
# of everything that needs to happen as a consequence:
#   - the data rows below it (old rows 3-9) shift up by one
#   - the shared-formula range on column AL is renumbered automatically
#   - the now-unused "Vestas_V47_public" shared string entry is dropped
#     and every other shared-string index used on the sheet is fixed up
#   - the sheet's dimension shrinks from A1:AT12 to A1:AT11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2").Delete()

# Excel leaves the selection sitting where the deleted row used to be,
# which after the shift corresponds to the (now mostly empty) row 13.
$ws.Range("B13").Select()
